# "valida[tion] and invalid testcase" workbook cleanup:
#  - rename sheet "NewCompanyDetails" -> "Sheet1" (via copy so the new
#    sheet gets a fresh sheetId, matching the authored edit)
#  - change header C1 from "res" to "exp" and give it a distinct font color
#  - drop the unused/blank D:F columns and the blank rows 7:12

$wb = $excel.ActiveWorkbook

$old = $wb.ActiveSheet
$oldName = $old.Name
$old.Copy($old) | Out-Null

$wb.Worksheets.Item($oldName).Delete() | Out-Null

$ws = $wb.Worksheets.Item("$oldName (2)")
$ws.Name = "Sheet1"

$ws.Range("C1").ClearFormats()
$ws.Range("C1").Value = "exp"
$ws.Range("C1").Font.Color = 4079210

$ws.Range("D1:F12").EntireColumn.Delete()
$ws.Range("A7:A12").EntireRow.Delete()
